# Cognitivity - The Backend: update the Introduction section text
# (ophir - updated backend report)

$d = $word.ActiveDocument

function Replace-RunText($Doc, $ParaIndex, $OldText, $NewText) {
    $para = $Doc.Paragraphs($ParaIndex)
    $found = $para.Range.Find.Execute($OldText, $true, $false, $false, $false, $false, `
                                       $true, 1, $false, $NewText, 2)
    if (-not $found) {
        throw "Replace-RunText: could not find text in paragraph $ParaIndex"
    }
}

function Split-ParagraphIntoTwoRuns($Doc, $ParaIndex, $Part1, $Part2) {
    $para = $Doc.Paragraphs($ParaIndex)
    $pRange = $para.Range
    $start = $pRange.Start
    $end = $pRange.End - 1   # exclude the paragraph mark

    # Replace the whole paragraph body (one run) with the full combined text;
    # this preserves the existing run formatting (sz/szCs/rtl).
    $whole = $Doc.Range($start, $end)
    $whole.Text = ($Part1 + $Part2)

    # Re-seat the end of the (now single) run after the text swap.
    $newEnd = $start + ($Part1 + $Part2).Length
    $splitPos = $start + $Part1.Length

    # Force a run break at $splitPos by nudging the font size of the 2nd half
    # away from and back to its original value - identical-format runs are
    # otherwise coalesced back into a single <w:r> by the writer.
    $secondHalf = $Doc.Range($splitPos, $newEnd)
    $origSize = $secondHalf.Font.Size
    $secondHalf.Font.Size = $origSize + 1
    $secondHalf.Font.Size = $origSize
}

# --- Paragraph 4: "The previous week it was decided ..." -> two runs ---
Split-ParagraphIntoTwoRuns $d 4 `
    "This following report reviews " `
    "the backend logic of our software."

# --- Paragraph 5: "As I was scoping the web ..." -> two runs + trailing empty run ---
Split-ParagraphIntoTwoRuns $d 5 `
    "The backend should b" `
    "e aligned with the reports of other team members, since some reports are on the frontend, defining an interface the backend will implement, and some are on finding the best data server to store our information. We need to find backend technologies that will:"

# Append an empty trailing run (rPr only has rtl) at the very end of paragraph 5,
# mirroring the empty runs used elsewhere in this document for paragraph padding.
$para5 = $d.Paragraphs(5)
$p5End = $para5.Range.End - 1
$tail = $d.Range($p5End, $p5End)
$tail.InsertAfter([string][char]0x200B)
$tailRange = $d.Range($p5End, $p5End + 1)
$tailRange.Font.Bold = 1
$tailRange.Font.Bold = 0
$tailRange.Text = ""

# --- Paragraph 6 (bullet item): "I need to find a backend technology ..." ---
Replace-RunText $d 6 `
    "I need to find a backend technology that is implemented in one of the allowed programming languages that fits the “backend/server side” logic, which has great features for development, and " `
    "Be implemented in one of the allowed programming languages that fits the “backend/server side” logic, have great features for development, and "

Replace-RunText $d 6 `
    " that Daniel will choose to be the best (for now, it seems that angularJS is the best pick for frontend)." `
    " that we will choose to be the best to use."

# --- Paragraph 7 (bullet item): "I need to find comfortable ways ..." ---
Replace-RunText $d 7 `
    "I need to find comfortable ways to use different types of libraries that help talking to a database. Since we will most likely be using some sql database, the programming language I’ll use for backend " `
    "Be able to talk to a database. Since we will most likely be using some sql database, the programming language we will use for backend "

Replace-RunText $d 7 `
    " supply easy enough ways to use this type of db. This should be verified with Guy, since he is working on deciding which type of data server we will use eventually." `
    " supply easy enough ways to use this type of db."

# --- Paragraph 9: "As the most commonly used language ..." ---
Replace-RunText $d 9 `
    " (among the allowed PL’s), this report will mainly focus on that language. Other language (such as C#) will be covered at the end of this report, if Daniel covers a technology that requires it (such as ASP)." `
    " (among the allowed PL’s), this report will mainly focus on that language."

Write-Output "done"
